$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Original data rows 35-39 (before this edit):
#   35: -487  Gorriti 3735
#   36: -488  Blanco Encalada 4896
#   37: 6465  AGUIRRE 368
#   38: -604  Gurruchaga 2126
#   39: 7359  VEGA, NICETO, CNEL. 4678
#
# Target result: only one row remains after row 34, containing the former
# row 36 (-488 Blanco Encalada) data, renumbered as row 35.
#
# Remove the trailing rows first (bottom-up) so row numbers of rows not yet
# deleted stay stable, then remove row 35 (Gorriti) so row 36 shifts up to
# become the new row 35.
$ws.Rows.Item(39).EntireRow.Delete()
$ws.Rows.Item(38).EntireRow.Delete()
$ws.Rows.Item(37).EntireRow.Delete()
$ws.Rows.Item(35).EntireRow.Delete()
